# Daily refresh of the "cryptos" price/volume snapshot (GitHub Actions job).
# Values below come straight from the new scrape; column D (Price) holds
# strings (it mixes thousand-dotted formats like "66.990.91" with plain
# decimals like "577.89"), so a leading "'" is used wherever the new text
# would otherwise be auto-parsed by Excel as a genuine number - this keeps
# the cell text-typed, matching the source data's original inline-string
# representation. Rows 38/39 (OKB vs. Stacks) also swap rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.990.91"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.121.23"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'577.89"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'173.57"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D9").Value = "'6.44"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "'37.26"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "3.639.67"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "66.947.90"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "3.120.84"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "'16.26"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "'475.74"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("E22").Value = "  +5.46%  "
$ws.Range("D23").Value = "'83.93"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "'13.30"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'7.93"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "'28.58"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "0.0₃0959"
$ws.Range("E33").Value = "  -6.69%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'5.84"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").Value = "'0.975"
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("D37").Value = "'47.20"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'50.21"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'2.06"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "2.816.21"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").Value = "'383.36"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("E46").Value = "  -8.86%  "
$ws.Range("D47").Value = "'135.38"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D49").Value = "'24.91"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("E51").Value = "  -0.58%  "
